$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Turn the single-line "Design." paragraph into four bulleted paragraphs
#    describing the PlayerThread / Card / CardDeck / CardGame design choices.
# ---------------------------------------------------------------------------

# Clear the existing text first (but keep the paragraph mark) so that when we
# restyle + retype it, Word doesn't drag the old 12pt direct formatting along
# with the new runs.
$designPara = $d.Paragraphs(2)
$designRange = $designPara.Range
$clearRange = $d.Range($designRange.Start, $designRange.End - 1)
$clearRange.Text = ""

$designPara.Style = "List Paragraph"
$r1 = $designPara.Range
$rr1 = $d.Range($r1.Start, $r1.End - 1)
$rr1.Text = 'PlayerThread class that implements Runnable, with each player running their own instructions for playing the game in parallel. This was done so that the program runs more efficiently and to simulate how the game would be played in real life, with each player thinking and acting independently. Each player stores a list of Card objects to represent the player''s hand, and it has methods to draw cards from the CardDeck object that represents the deck to the players left, and similarly to discard cards to the deck on the right. The references to these decks are given to the player when it is instantiated.'
$rr1.InsertParagraphAfter()

$p2 = $d.Paragraphs(3)
$p2.Style = "List Paragraph"
$r2 = $p2.Range
$rr2 = $d.Range($r2.Start, $r2.End - 1)
$rr2.Text = 'Card class with each card containing an attribute to store its numerical value. This was done instead of just using integers so that multiple cards with the same value are still unique objects.'
$rr2.InsertParagraphAfter()

$p3 = $d.Paragraphs(4)
$p3.Style = "List Paragraph"
$r3 = $p3.Range
$rr3 = $d.Range($r3.Start, $r3.End - 1)
$rr3.Text = 'CardDeck class representing a deck with four cards which are changed throughout the course of the game. Each deck stores its cards in a queue of Card objects because cards are always added to a deck at the bottom and always taken from the top, giving it a first-in-first-out rule of operations.'
$rr3.InsertParagraphAfter()

$p4 = $d.Paragraphs(5)
$p4.Style = "List Paragraph"
$r4 = $p4.Range
$rr4 = $d.Range($r4.Start, $r4.End - 1)
$rr4.Text = 'CardGame class where the game is run, containing the main() method. The game creates a list of PlayerThread objects and CardDeck objects that each have a specified length, taken from a user input. These PlayerThread threads are run so that each player plays the game in parallel, at which point this main() method stops performing any computation until a player has won. When the game has finished this method will output which player has won the game and the program will stop.'

# Turn the whole four-paragraph block into one single bulleted list so they
# all share the same numbering definition (numId).
$listRange = $d.Range($d.Paragraphs(2).Range.Start, $d.Paragraphs(5).Range.End)
$listRange.ListFormat.ApplyBulletDefault()

# Tune the auto-created "List Paragraph" style so it matches Word's built-in
# definition (left indent + contextual spacing, ui priority 34).
$listStyle = $d.Styles("List Paragraph")
$listStyle.Priority = 34
$listStyle.ParagraphFormat.LeftIndent = 36
$listStyle.NoSpaceBetweenParagraphsOfSameStyle = $true

# ---------------------------------------------------------------------------
# 2. Drop the explicit 12pt (sz=24) direct formatting from the "Test." line
#    so it just inherits the Normal style like the rest of the doc.
# ---------------------------------------------------------------------------

$testPara = $d.Paragraphs($d.Paragraphs.Count)
$testRange = $testPara.Range
$testClear = $d.Range($testRange.Start, $testRange.End - 1)
$testClear.Text = ""
$testPara.Style = "Normal"
$testRange2 = $testPara.Range
$testFinal = $d.Range($testRange2.Start, $testRange2.End - 1)
$testFinal.Text = "Test."

Write-Host "done"
